$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$hf = $m.HeadersFooters
$hf.Footer.Text = "<footer>"
$hf.SlideNumber.Text = "<number>"
$hf.DateAndTime.Text = "<date/time>"
Write-Host "Footer=[$($hf.Footer.Text)]"
Write-Host "SlideNumber=[$($hf.SlideNumber.Text)]"
Write-Host "DateAndTime=[$($hf.DateAndTime.Text)]"
